# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to match the freshly scraped data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 5,6,7,9,10,11 -> column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 161
$wsExpo.Range("F6").Value = 170
$wsExpo.Range("F7").Value = 310
$wsExpo.Range("F9").Value = 2114
$wsExpo.Range("F10").Value = 366
$wsExpo.Range("F11").Value = 5110

# Sheet "全部类型": rows 6,7,8,12,13,14 -> column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 161
$wsAll.Range("F7").Value = 170
$wsAll.Range("F8").Value = 310
$wsAll.Range("F12").Value = 2114
$wsAll.Range("F13").Value = 366
$wsAll.Range("F14").Value = 5110
